$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AHB-Diff")

# --- Rename header cells: "_old" -> "_FV2310", "_new" -> "_FV2404" ---
$oldSuffixCols = @("A","B","C","D","E","F","G","H","I","J")
$newSuffixCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2310")
}

foreach ($col in $newSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2404")
}

# --- Turn the used range into an Excel Table (ListObject) ---
$rng = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- Freeze the header row ---
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.FreezePanes = $true
